# OSAT.xlsx: add a new analysis row (row 5) by duplicating row 4 ("PI082_02")
# down one row - same sample data, same style, formulas re-pointed at row 5 -
# then move the selection to BJ7 (fix bugs: excel AND statm & Fo calc, after Gavrilenko).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns A:V on row 4 hold literal sample values (sample_no + oxide wt% + totals).
$valueCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

# Columns W:BZ on row 4 hold the derived formulas (cation fractions, Fo, thermometry, Z-test, ...).
$formulaCols = @("W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD","BE","BF","BG","BH","BI","BJ","BK","BL","BM","BN","BO","BP","BQ","BR","BS","BT","BU","BV","BW","BX","BY","BZ")

foreach ($col in $valueCols) {
    $srcCell = $ws.Range($col + "4")
    $dstCell = $ws.Range($col + "5")
    # Value2 (not Text, which is display-rounded) to keep full double precision.
    $dstCell.Value = $srcCell.Value2
}

foreach ($col in $formulaCols) {
    $srcCell = $ws.Range($col + "4")
    $dstCell = $ws.Range($col + "5")
    $srcFormula = $srcCell.Formula
    # same-row-relative refs only move from row 4 -> row 5 (Sheet2! anchors stay put)
    $dstFormula = $srcFormula -replace '([A-Z]{1,3})4', '${1}5'
    $dstCell.Formula = $dstFormula
}

# Land the selection where the user ended up after adding the row.
$ws.Activate()
$ws.Range("BJ7").Select()
